$d = $word.ActiveDocument

# --- 1. Table cell margins: left margin 30 dxa (1.5 pt) -> 27 dxa (1.35 pt) ---
# 30 dxa / 20 = 1.5 pt ; 27 dxa / 20 = 1.35 pt
$t = $d.Tables.Item(1)
$t.LeftPadding = 1.35
foreach ($row in $t.Rows) {
  foreach ($cell in $row.Cells) {
    $cell.LeftPadding = 1.35
  }
}

# --- 2. Study 6 results paragraph: correct the regression coefficients ---
# "b = 1.16" was a typo; corrected value is "b = 1.74"
$d.Content.Find.Execute("= 1.16, ", $true, $false, $false, $false, $false, $true, 1, $false, "= 1.74, ", 2)

# "SE = 1.66" was a typo; corrected value is "SE = 2.49"
$d.Content.Find.Execute("= 1.66, ", $true, $false, $false, $false, $false, $true, 1, $false, "= 2.49, ", 2)
